# Auto-generated edit script applying numeric corrections to leve profit sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = $null
$ws.Range("N43").Value = -3138
$ws.Range("H125").Value = 2872.6667
$ws.Range("J125").Value = 571.4286
$ws.Range("L125").Value = 5142.8574
$ws.Range("N125").Value = -10062.8574
$ws.Range("H127").Value = 1475.0233
$ws.Range("J127").Value = 2010.2858
$ws.Range("L127").Value = 6030.857400000001
$ws.Range("N127").Value = -15950.8574
$ws.Range("H134").Value = 31000
$ws.Range("J134").Value = 31000
$ws.Range("L134").Value = 31000
$ws.Range("N134").Value = -41140
$ws.Range("H136").Value = 25250
$ws.Range("J136").Value = 25250
$ws.Range("L136").Value = 25250
$ws.Range("N136").Value = -35450
$ws.Range("H137").Value = 1711.1708
$ws.Range("I137").Value = 1161.8572
$ws.Range("J137").Value = 1996
$ws.Range("K137").Value = 3485.5716
$ws.Range("L137").Value = 5988
$ws.Range("M137").Value = -935.5715999999998
$ws.Range("N137").Value = -11088

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9446.056
$ws.Range("I45").Value = 10955.934
$ws.Range("K45").Value = 10955.934
$ws.Range("M45").Value = -10578.934
$ws.Range("H61").Value = 9241.25
$ws.Range("I61").Value = 3275.7144
$ws.Range("J61").Value = 51000
$ws.Range("K61").Value = 3275.7144
$ws.Range("L61").Value = 51000
$ws.Range("M61").Value = -3063.7144
$ws.Range("N61").Value = -51424
$ws.Range("H92").Value = 37516.668
$ws.Range("J92").Value = 37516.668
$ws.Range("L92").Value = 37516.668
$ws.Range("N92").Value = -42508.668
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null
$ws.Range("H133").Value = 39326.832
$ws.Range("J133").Value = 39326.832
$ws.Range("L133").Value = 39326.832
$ws.Range("N133").Value = -44386.832
$ws.Range("H136").Value = 9241.25
$ws.Range("I136").Value = 3275.7144
$ws.Range("J136").Value = 51000
$ws.Range("K136").Value = 9827.143199999999
$ws.Range("L136").Value = 153000
$ws.Range("M136").Value = -7277.143199999999
$ws.Range("N136").Value = -158100
$ws.Range("H139").Value = 37844.875
$ws.Range("J139").Value = 37844.875
$ws.Range("L139").Value = 37844.875
$ws.Range("N139").Value = -48124.875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 12481.685
$ws.Range("I105").Value = 18359
$ws.Range("J105").Value = 2406.2856
$ws.Range("K105").Value = 18359
$ws.Range("L105").Value = 2406.2856
$ws.Range("M105").Value = -16612
$ws.Range("N105").Value = -5900.2856
$ws.Range("H107").Value = 1342.8572
$ws.Range("I107").Value = 1342.8572
$ws.Range("K107").Value = 1342.8572
$ws.Range("M107").Value = 577.1428000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1986.1163
$ws.Range("I134").Value = 2077
$ws.Range("J134").Value = 1100
$ws.Range("K134").Value = 6231
$ws.Range("L134").Value = 3300
$ws.Range("M134").Value = -3696
$ws.Range("N134").Value = -8370

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7225.975
$ws.Range("I3").Value = 9574.214
$ws.Range("J3").Value = 5961.5386
$ws.Range("K3").Value = 28722.642
$ws.Range("L3").Value = 17884.6158
$ws.Range("M3").Value = -28610.642
$ws.Range("N3").Value = -18108.6158
$ws.Range("H5").Value = 139481.77
$ws.Range("J5").Value = 172391.62
$ws.Range("L5").Value = 517174.86
$ws.Range("N5").Value = -517398.86
$ws.Range("H68").Value = 1870.9656
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1870.9656
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 5612.8968
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = -7234.8968
$ws.Range("H71").Value = 1870.9656
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1870.9656
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 16838.6904
$ws.Range("M71").Value = $null
$ws.Range("N71").Value = -24950.6904
$ws.Range("H131").Value = 18966602
$ws.Range("J131").Value = 23810844
$ws.Range("L131").Value = 71432532
$ws.Range("N131").Value = -71442612
$ws.Range("H134").Value = 9279.413
$ws.Range("I134").Value = 7928.3125
$ws.Range("K134").Value = 23784.9375
$ws.Range("M134").Value = -18714.9375
$ws.Range("H135").Value = 139481.77
$ws.Range("J135").Value = 172391.62
$ws.Range("L135").Value = 1551524.58
$ws.Range("N135").Value = -1556594.58
$ws.Range("H136").Value = 15894.286
$ws.Range("I136").Value = 17376.666
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 52129.99800000001
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -47029.99800000001
$ws.Range("N136").Value = -31200
$ws.Range("H137").Value = 27040772
$ws.Range("I137").Value = 1168.4166
$ws.Range("J137").Value = 40019780
$ws.Range("K137").Value = 3505.2498
$ws.Range("L137").Value = 120059340
$ws.Range("M137").Value = 1594.7502
$ws.Range("N137").Value = -120069540
$ws.Range("H138").Value = 7363
$ws.Range("I138").Value = 7486.6665
$ws.Range("J138").Value = 6250
$ws.Range("K138").Value = 22459.9995
$ws.Range("L138").Value = 18750
$ws.Range("M138").Value = -17319.9995
$ws.Range("N138").Value = -29030
$ws.Range("H139").Value = 34229.707
$ws.Range("I139").Value = 69803.125
$ws.Range("J139").Value = 2608.889
$ws.Range("K139").Value = 209409.375
$ws.Range("L139").Value = 7826.667
$ws.Range("M139").Value = -204269.375
$ws.Range("N139").Value = -18106.667

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4594263.5
$ws.Range("I122").Value = 2819763.5
$ws.Range("J122").Value = 25001012
$ws.Range("K122").Value = 8459290.5
$ws.Range("L122").Value = 75003036
$ws.Range("M122").Value = -8456840.5
$ws.Range("N122").Value = -75007936
$ws.Range("H132").Value = 3784.5925
$ws.Range("I132").Value = 3392.4707
$ws.Range("J132").Value = 4451.2
$ws.Range("K132").Value = 10177.4121
$ws.Range("L132").Value = 13353.6
$ws.Range("M132").Value = -7647.4121
$ws.Range("N132").Value = -18413.6

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5678.577
$ws.Range("I7").Value = 2479.889
$ws.Range("J7").Value = 12875.625
$ws.Range("K7").Value = 2479.889
$ws.Range("L7").Value = 12875.625
$ws.Range("M7").Value = -2367.889
$ws.Range("N7").Value = -13099.625
$ws.Range("H40").Value = 500002240
$ws.Range("I40").Value = 500002240
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 500002240
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -500002104
$ws.Range("N40").Value = $null
$ws.Range("H126").Value = 5678.577
$ws.Range("I126").Value = 2479.889
$ws.Range("J126").Value = 12875.625
$ws.Range("K126").Value = 7439.667
$ws.Range("L126").Value = 38626.875
$ws.Range("M126").Value = -4969.667
$ws.Range("N126").Value = -43566.875
$ws.Range("H133").Value = 79403.71000000001
$ws.Range("J133").Value = 79403.71000000001
$ws.Range("L133").Value = 79403.71000000001
$ws.Range("N133").Value = -84463.71000000001
$ws.Range("H136").Value = 5517.6763
$ws.Range("I136").Value = 2739.48
$ws.Range("J136").Value = 13234.889
$ws.Range("K136").Value = 8218.440000000001
$ws.Range("L136").Value = 39704.667
$ws.Range("M136").Value = -5668.440000000001
$ws.Range("N136").Value = -44804.667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 250007950
$ws.Range("I107").Value = 333333920
$ws.Range("J107").Value = 30003
$ws.Range("K107").Value = 1000001760
$ws.Range("L107").Value = 90009
$ws.Range("M107").Value = -999999840
$ws.Range("N107").Value = -93849
